$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns after AB1, matching the existing header style
$ws.Range("AC1").Value = "Diad2_prom/std_betweendiads"
$ws.Range("AD1").Value = "Left_vs_Right"

# Match the formatting of the last existing header cell (AB1)
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AD1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
